$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.924.95"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.448.63"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "2.448.69"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -7.58%  "
$ws.Range("D14").Value = "2.884.04"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "57.551.42"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "2.450.46"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.67%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "2.545.75"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.156"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("B29").Value = "Polygon"
$ws.Range("C29").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.381"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.17%  "
$ws.Range("E30").Value = "  +4.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "0.0₃0733"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("E39").Value = "  +4.46%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.39%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "255.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.572"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0917"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0491"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  +1.13%  "

